$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.355.14"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -4.32%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.759.91"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -4.07%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9999"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "303.96"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.61%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4281"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.22%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3619"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.87%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07054"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.76%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8310"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.73%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.16"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.30%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.776.15"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.61%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.217"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.81%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.354"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.38%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06786"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.46%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "79.13"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.65%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008629"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.09%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.00%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.91"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.255.34"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.71%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.980"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.47%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.05"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.97%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.977.92"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.68%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.907"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.24%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.87"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.90%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.08"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.26%  "

$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.032"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.03%  "

$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.50"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.21%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.681"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -7.75%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08869"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.28%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7236"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.15%  "

$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.118"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.94%  "

$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.306"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.90%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.000"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.711"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -9.24%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.068"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.10%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05100"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.94%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01878"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.88%  "

$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4898"

$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1601"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.17%  "

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.168"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.81%  "

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.461"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -12.07%  "

$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.017"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.61%  "

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.66"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.68%  "

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9999"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.02%  "

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.966"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.29%  "

$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06178"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.64%  "

$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Decentraland"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4471"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.41%  "

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.569"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.94%  "

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.716"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.10%  "
